# Updates "Horarios 141" workbook with the latest scrape (06:58:58)
# Sheet 1: LP1912, Sheet 2: LP1912-215, Sheet 3: 6203-6173

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: 06:58:58"
$ws1.Range("A3").Value = "Total filas: 45"

$rows1 = @(
    @(15, "06:58:58", "06:58", "215A_EL PATO", 0, "LP1912"),
    @(16, "06:58:58", "06:58", "225_GOMEZ", 0, "LP1912"),
    @(17, "06:46:40", "06:59", "225_GOMEZ", 13, "LP1912"),
    @(18, "06:58:58", "07:15", "215C_EL PATO", 17, "LP1912"),
    @(19, "06:15:23", "07:16", "215C_EL PATO", 61, "LP1912"),
    @(20, "06:58:58", "07:18", "14_ABASTO", 20, "LP1912"),
    @(21, "06:46:40", "07:19", "14_ABASTO", 33, "LP1912"),
    @(22, "06:58:58", "07:20", "16_SANTA ANA", 22, "LP1912"),
    @(23, "06:15:23", "07:21", "16_SANTA ANA", 66, "LP1912"),
    @(24, "06:46:40", "07:21", "23_HERNANDEZ", 35, "LP1912"),
    @(25, "06:58:58", "07:22", "23_HERNANDEZ", 24, "LP1912"),
    @(26, "06:58:58", "07:29", "17X38_ROMERO", 31, "LP1912"),
    @(27, "06:58:58", "07:34", "10_OLMOS", 36, "LP1912"),
    @(28, "06:46:40", "07:35", "10_OLMOS", 49, "LP1912"),
    @(29, "06:58:58", "07:36", "27_EL RETIRO", 38, "LP1912"),
    @(30, "06:15:23", "07:37", "27_EL RETIRO", 82, "LP1912"),
    @(31, "06:58:58", "07:43", "215A_EL PATO", 45, "LP1912"),
    @(32, "06:58:58", "07:54", "14_ABASTO", 56, "LP1912"),
    @(33, "06:46:40", "07:55", "14_ABASTO", 69, "LP1912"),
    @(34, "06:58:58", "07:59", "17_ROMERO", 61, "LP1912"),
    @(35, "06:46:40", "08:00", "17_ROMERO", 74, "LP1912"),
    @(36, "06:58:58", "08:00", "16_SANTA ANA", 62, "LP1912"),
    @(37, "06:46:40", "08:01", "16_SANTA ANA", 75, "LP1912"),
    @(38, "06:46:40", "08:06", "23_HERNANDEZ", 80, "LP1912"),
    @(39, "06:58:58", "08:11", "10_OLMOS", 73, "LP1912"),
    @(40, "06:58:58", "08:12", "15X38_ABASTO", 74, "LP1912"),
    @(41, "06:46:40", "08:13", "15X38_ABASTO", 87, "LP1912"),
    @(42, "06:58:58", "08:28", "15_ABASTO", 90, "LP1912"),
    @(43, "06:58:58", "08:28", "11_ETCHEVERRY", 90, "LP1912"),
    @(44, "06:46:40", "08:29", "11_ETCHEVERRY", 103, "LP1912"),
    @(45, "06:46:40", "08:29", "15_ABASTO", 103, "LP1912"),
    @(46, "06:58:58", "08:40", "16_P MOR-SANTA ANA", 102, "LP1912"),
    @(47, "06:46:40", "08:41", "16_P MOR-SANTA ANA", 115, "LP1912"),
    @(48, "06:58:58", "08:43", "215C_EL PATO", 105, "LP1912"),
    @(49, "06:58:58", "08:52", "23_HERNANDEZ", 114, "LP1912"),
    @(50, "06:58:58", "08:53", "215B_EL PATO", 115, "LP1912")
)

foreach ($row in $rows1) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: 06:58:58"
$ws2.Range("A3").Value = "Total filas: 9"

$rows2 = @(
    @(9, "06:58:58", "06:58", "215A_EL PATO", 0, "LP1912"),
    @(10, "06:58:58", "07:15", "215C_EL PATO", 17, "LP1912"),
    @(11, "06:15:23", "07:16", "215C_EL PATO", 61, "LP1912"),
    @(12, "06:58:58", "07:43", "215A_EL PATO", 45, "LP1912"),
    @(13, "06:58:58", "08:43", "215C_EL PATO", 105, "LP1912"),
    @(14, "06:58:58", "08:53", "215B_EL PATO", 115, "LP1912")
)

foreach ($row in $rows2) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 3).Value = $row[3]
    $ws2.Cells.Item($r, 4).Value = $row[4]
    $ws2.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: 06:58:58"
$ws3.Range("A3").Value = "Total filas: 4"

$rows3 = @(
    @(6, "06:58:58", "07:42", "215A_LA PLATA", 44, "L6173"),
    @(8, "06:58:58", "08:35", "215A_LA PLATA", 97, "L6173"),
    @(9, "06:58:58", "08:50", "215C_LA PLATA", 112, "L6203")
)

foreach ($row in $rows3) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "done"
